$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 427 (shifts existing rows 427-476 down to 428-477)
$ws.Rows.Item(427).Insert()

# Populate the newly inserted row 427 with the new weekly price record
$ws.Cells.Item(427, 1).Value = 3
$ws.Cells.Item(427, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(427, 3).Value = "Coquimbo"
$ws.Cells.Item(427, 4).Value = 44946
$ws.Cells.Item(427, 5).Value = 5
$ws.Cells.Item(427, 6).Value = 100112009
$ws.Cells.Item(427, 7).Value = "Acelga"
$ws.Cells.Item(427, 8).Value = "Sin especificar"
$ws.Cells.Item(427, 9).Value = "Primera"
$ws.Cells.Item(427, 10).Value = 220
$ws.Cells.Item(427, 11).Value = 3300
$ws.Cells.Item(427, 12).Value = 3500
$ws.Cells.Item(427, 13).Value = 3400
$ws.Cells.Item(427, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(427, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(427, 16).Value = 567
$ws.Cells.Item(427, 17).Value = 6
$ws.Cells.Item(427, 18).Value = "Hortaliza"
